$d = $word.ActiveDocument

function Replace-Text($old, $new) {
    $d.Content.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2) | Out-Null
}

# 1. "List with error with 3 items (...)" paragraph - rewritten wording.
Replace-Text 'List with error with 3 items (here the same example as previous one, you decribe list as optional in your data structure BUT you describe first item as required in your data structure resulting an error if you skip to insert it in your data):' 'List full optional with error with 3 items (since you decribe list as optional in your data structure BUT you describe first item as required in your data structure, it will throw an error if you skip to insert the first item in your data):'

# 2. Drop the leading "- " bullet-style prefix on the "Nom" line.
Replace-Text '- Nom : {$person.fields.name}, Age : {$person.fields.age}' 'Nom : {$person.fields.name}, Age : {$person.fields.age}'

# 3. Drop the leading "- " bullet-style prefix on the "Niveau de menace" line.
Replace-Text '- Niveau de menace : {$security.fields.security_lvl}, code menace : {$security.fields.code}' 'Niveau de menace : {$security.fields.security_lvl}, code menace : {$security.fields.code}'

# 4. "List with optional object ... with error (...)" paragraph - rewritten wording
#    (also normalizes the non-breaking space before "with 3 items" to a regular space).
Replace-Text ([char]0x00A0) ' '
Replace-Text 'List with optional object with 3 items with error (here the same example as previous one, you decribe list as required in your data structure BUT you describe last object as optional and one property of this object as required in your data structure resulting an error if you skip to insert it in your data):' 'List with optional object with 3 items with error (since you decribe list as required in your data structure BUT you describe last object of the list as optional and one property of this object as required in your data structure, it will throw an error if you skip to insert this field in your data):'

# 5. Drop the leading "- " bullet-style prefix on the "Pays" line.
Replace-Text '- Pays : {$destination.fields.country}, Etat: {$destination.fields.state}' 'Pays : {$destination.fields.country}, Etat: {$destination.fields.state}'

# 6. Drop the leading "- " bullet-style prefix on the nested "{$person.fields.name}" line.
Replace-Text '- {$person.fields.name}' '{$person.fields.name}'

# 7. Drop the leading "- " bullet-style prefix on the nested "{$task}" line.
Replace-Text '- {$task}' '{$task}'

# 8. The paragraph holding "{END-FOR person}" gains explicit (false) bold / bold-complex-script
#    paragraph-mark run properties (<w:rPr><w:b w:val="false"/><w:bCs w:val="false"/></w:rPr>
#    under <w:pPr>), while the run text itself keeps its default (inherited) formatting.
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.TrimEnd([char]13, [char]7) -eq '{END-FOR person}') {
        $full = $p.Range
        $full.Font.Bold = $false
        $full.Font.BoldBi = $false
        # Re-clear the formatting on the visible text (everything but the trailing
        # paragraph mark) so only the paragraph mark itself keeps the explicit
        # "not bold" flag - leaves the run's own <w:rPr/> untouched/empty.
        $textOnly = $d.Range($full.Start, $full.End - 1)
        $textOnly.Font.Bold = 9999999
    }
}
